# Applies the cryptos-list price/volume refresh described in the commit
# "Updated cryptos list on Sun Feb 18 13:49:47 UTC 2024 with GitHub Actions".
#
# Note: column D values are stored as literal text (e.g. "51.616.85", "352.37")
# even though many of them look like numbers. Assigning a plain numeric-looking
# string via .Value lets Excel auto-convert it to a real number (e.g. "26.00" ->
# 26), which would lose the original text formatting. To keep such values as text
# (matching the workbook author's intent), we prefix them with a leading
# apostrophe, exactly like typing them into Excel by hand; Excel strips the
# apostrophe and stores the remainder as a text cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.612.43'
$ws.Range('E2').Value = '  +0.69%  '

$ws.Range('D3').Value = '2.791.54'
$ws.Range('E3').Value = '  +1.57%  '

$ws.Range('D5').Value = '''352.30'
$ws.Range('E5').Value = '  -0.82%  '

$ws.Range('D6').Value = '''111.13'
$ws.Range('E6').Value = '  +2.94%  '

$ws.Range('E7').Value = '  +0.69%  '

$ws.Range('E8').Value = '  -0.09%  '

$ws.Range('E9').Value = '  +8.41%  '

$ws.Range('D10').Value = '''39.92'
$ws.Range('E10').Value = '  +1.92%  '

$ws.Range('E11').Value = '  -1.10%  '

$ws.Range('E12').Value = '  +0.07%  '

$ws.Range('E13').Value = '  +1.67%  '

$ws.Range('E14').Value = '  +3.33%  '

$ws.Range('D15').Value = '3.231.99'
$ws.Range('E15').Value = '  +1.45%  '

$ws.Range('D16').Value = '2.792.40'
$ws.Range('E16').Value = '  +0.46%  '

$ws.Range('E17').Value = '  +2.17%  '

$ws.Range('D18').Value = '51.596.57'
$ws.Range('E18').Value = '  +0.81%  '

$ws.Range('E19').Value = '  +0.72%  '

$ws.Range('D20').Value = '''3.19'
$ws.Range('E20').Value = '  +6.31%  '

$ws.Range('D21').Value = '''13.51'
$ws.Range('E21').Value = '  +3.92%  '

$ws.Range('D22').Value = '0.0₃0969'
$ws.Range('E22').Value = '  +1.19%  '

$ws.Range('D23').Value = '''70.13'
$ws.Range('E23').Value = '  +1.00%  '

$ws.Range('D24').Value = '''266.95'
$ws.Range('E24').Value = '  +0.81%  '

$ws.Range('E25').Value = '  +0.13%  '

$ws.Range('D26').Value = '''0.999'
$ws.Range('E26').Value = '  -0.06%  '

$ws.Range('D27').Value = '''26.00'
$ws.Range('E27').Value = '  +0.15%  '

$ws.Range('E28').Value = '  -1.27%  '

$ws.Range('D29').Value = '''38.86'
$ws.Range('E29').Value = '  +11.31%  '

$ws.Range('D30').Value = '''10.32'
$ws.Range('E30').Value = '  +2.58%  '

$ws.Range('E31').Value = '  +0.47%  '

$ws.Range('D32').Value = '''52.58'
$ws.Range('E32').Value = '  +2.35%  '

$ws.Range('D33').Value = '''6.10'
$ws.Range('E33').Value = '  +0.64%  '

$ws.Range('D34').Value = '''0.0454'
$ws.Range('E34').Value = '  +3.10%  '

$ws.Range('D35').Value = '''0.0886'
$ws.Range('E35').Value = '  +6.83%  '

$ws.Range('D36').Value = '''5.53'
$ws.Range('E36').Value = '  +7.67%  '

$ws.Range('E37').Value = '  +0.00%  '

$ws.Range('D38').Value = '''18.72'

$ws.Range('E39').Value = '  +3.29%  '

$ws.Range('E40').Value = '  +0.97%  '

$ws.Range('E41').Value = '  +1.36%  '

$ws.Range('E42').Value = '  +0.11%  '

$ws.Range('E43').Value = '  +0.78%  '

$ws.Range('D44').Value = '''120.01'
$ws.Range('E44').Value = '  +0.66%  '

$ws.Range('D45').Value = '''21.79'
$ws.Range('E45').Value = '  -0.77%  '

$ws.Range('B46').Value = 'NEARProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D46').Value = '''3.40'
$ws.Range('E46').Value = '  +5.63%  '

$ws.Range('B47').Value = 'ApeXProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D47').Value = '''2.43'
$ws.Range('E47').Value = '  +5.77%  '

$ws.Range('D48').Value = '2.101.67'
$ws.Range('E48').Value = '  +1.10%  '

$ws.Range('D49').Value = '''0.952'
$ws.Range('E49').Value = '  +3.10%  '

$ws.Range('E50').Value = '  -0.29%  '

$ws.Range('E51').Value = '  +6.42%  '

Write-Host "Applied cryptos update"